$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,2).Value = "incorrect"
$ws.Cells.Item(2,3).Value = "correct"
$ws.Cells.Item(2,4).Value = 3.835552596714483
$ws.Cells.Item(2,5).Value = 0.8900544242691064
$ws.Cells.Item(2,6).Value = 82.51391365710614
$ws.Cells.Item(2,7).Value = 21.11093335612109

$ws.Cells.Item(3,2).Value = "incorrect"
$ws.Cells.Item(3,3).Value = "correct"
$ws.Cells.Item(3,4).Value = 5.528416835309878
$ws.Cells.Item(3,5).Value = 0.9086528341696638
$ws.Cells.Item(3,6).Value = 66.51081684824972
$ws.Cells.Item(3,7).Value = 32.5356607209018

$ws.Cells.Item(4,2).Value = "incorrect"
$ws.Cells.Item(4,3).Value = "correct"
$ws.Cells.Item(4,4).Value = 3.389014367264778
$ws.Cells.Item(4,5).Value = 4.125642136441833
$ws.Cells.Item(4,6).Value = 57.59293919951686
$ws.Cells.Item(4,7).Value = 4.003238373410056

$ws.Cells.Item(5,2).Value = "incorrect"
$ws.Cells.Item(5,3).Value = "correct"
$ws.Cells.Item(5,4).Value = 1.14708997008554
$ws.Cells.Item(5,5).Value = 3.076693961165187
$ws.Cells.Item(5,6).Value = 32.82583326047829
$ws.Cells.Item(5,7).Value = 82.70594722083857

$ws.Cells.Item(6,2).Value = "incorrect"
$ws.Cells.Item(6,3).Value = "correct"
$ws.Cells.Item(6,4).Value = 5.472265686199302
$ws.Cells.Item(6,5).Value = 2.66442582073238
$ws.Cells.Item(6,6).Value = 46.90064166430977
$ws.Cells.Item(6,7).Value = 27.85683433058443

$ws.Cells.Item(7,2).Value = "incorrect"
$ws.Cells.Item(7,3).Value = "correct"
$ws.Cells.Item(7,4).Value = 1.334743394907734
$ws.Cells.Item(7,5).Value = 3.896740308473223
$ws.Cells.Item(7,6).Value = 90.00417312256316
$ws.Cells.Item(7,7).Value = 33.85122998272168

$ws.Cells.Item(8,2).Value = "incorrect"
$ws.Cells.Item(8,3).Value = "correct"
$ws.Cells.Item(8,4).Value = 5.60706370553178
$ws.Cells.Item(8,5).Value = 3.981403706453577
$ws.Cells.Item(8,6).Value = 69.17012194640326
$ws.Cells.Item(8,7).Value = 44.32275577968915

$ws.Cells.Item(9,2).Value = "incorrect"
$ws.Cells.Item(9,3).Value = "incorrect"
$ws.Cells.Item(9,4).Value = 2.658308209248342
$ws.Cells.Item(9,5).Value = 4.869141274352523
$ws.Cells.Item(9,6).Value = 109.0599159400299
$ws.Cells.Item(9,7).Value = 44.79017057808762

$ws.Cells.Item(10,2).Value = "incorrect"
$ws.Cells.Item(10,3).Value = "correct"
$ws.Cells.Item(10,4).Value = 3.791839669140422
$ws.Cells.Item(10,5).Value = 3.198720736692172
$ws.Cells.Item(10,6).Value = 59.55194981035638
$ws.Cells.Item(10,7).Value = 17.85302097753112

$ws.Cells.Item(11,2).Value = "correct"
$ws.Cells.Item(11,3).Value = "incorrect"
$ws.Cells.Item(11,4).Value = 4.568662008538956
$ws.Cells.Item(11,5).Value = 2.512084667368572
$ws.Cells.Item(11,6).Value = 97.52252416630938
$ws.Cells.Item(11,7).Value = 15.09202414934451

$ws.Cells.Item(12,2).Value = "incorrect"
$ws.Cells.Item(12,3).Value = "correct"
$ws.Cells.Item(12,4).Value = 3.884035779029048
$ws.Cells.Item(12,5).Value = 1.786269840004332
$ws.Cells.Item(12,6).Value = 67.80077229457494
$ws.Cells.Item(12,7).Value = 29.98694873542614

$ws.Cells.Item(13,2).Value = "incorrect"
$ws.Cells.Item(13,3).Value = "correct"
$ws.Cells.Item(13,4).Value = 3.169734189657093
$ws.Cells.Item(13,5).Value = 1.252430595416194
$ws.Cells.Item(13,6).Value = 57.48898754577474
$ws.Cells.Item(13,7).Value = 38.87121130161695

$ws.Cells.Item(14,2).Value = "incorrect"
$ws.Cells.Item(14,3).Value = "incorrect"
$ws.Cells.Item(14,4).Value = 2.654888525402939
$ws.Cells.Item(14,5).Value = 2.595486170518498
$ws.Cells.Item(14,6).Value = 111.0655183697489
$ws.Cells.Item(14,7).Value = 17.25448013471419

$ws.Cells.Item(15,2).Value = "incorrect"
$ws.Cells.Item(15,3).Value = "correct"
$ws.Cells.Item(15,4).Value = 4.987705425647089
$ws.Cells.Item(15,5).Value = 3.01075022968644
$ws.Cells.Item(15,6).Value = 78.40780586747414
$ws.Cells.Item(15,7).Value = 9.622918861120127

$ws.Cells.Item(16,2).Value = "correct"
$ws.Cells.Item(16,3).Value = "correct"
$ws.Cells.Item(16,4).Value = 4.040161468576324
$ws.Cells.Item(16,5).Value = 3.177873423813282
$ws.Cells.Item(16,6).Value = 47.50370270280273
$ws.Cells.Item(16,7).Value = 23.25523237579559

$ws.Cells.Item(17,2).Value = "correct"
$ws.Cells.Item(17,3).Value = "correct"
$ws.Cells.Item(17,4).Value = 4.611214123034634
$ws.Cells.Item(17,5).Value = 1.880188157187829
$ws.Cells.Item(17,6).Value = 14.15151714614174
$ws.Cells.Item(17,7).Value = 60.97035371558065

$ws.Cells.Item(18,2).Value = "incorrect"
$ws.Cells.Item(18,3).Value = "correct"
$ws.Cells.Item(18,4).Value = 5.202128301578574
$ws.Cells.Item(18,5).Value = 1.708370936279019
$ws.Cells.Item(18,6).Value = 79.49208813250146
$ws.Cells.Item(18,7).Value = 24.4913039108656

$ws.Cells.Item(19,2).Value = "correct"
$ws.Cells.Item(19,3).Value = "incorrect"
$ws.Cells.Item(19,4).Value = 2.40317936485656
$ws.Cells.Item(19,5).Value = 3.028417202288276
$ws.Cells.Item(19,6).Value = 67.95614306222794
$ws.Cells.Item(19,7).Value = 49.12876283386283

$ws.Cells.Item(20,2).Value = "correct"
$ws.Cells.Item(20,3).Value = "correct"
$ws.Cells.Item(20,4).Value = 5.773410114174354
$ws.Cells.Item(20,5).Value = 1.145593332901836
$ws.Cells.Item(20,6).Value = 30.92956563070015
$ws.Cells.Item(20,7).Value = 52.32717916647265

$ws.Cells.Item(21,2).Value = "incorrect"
$ws.Cells.Item(21,3).Value = "correct"
$ws.Cells.Item(21,4).Value = 2.483300150669507
$ws.Cells.Item(21,5).Value = 3.119946118303702
$ws.Cells.Item(21,6).Value = 73.03745716877749
$ws.Cells.Item(21,7).Value = 2.36392095913525

$ws.Cells.Item(22,2).Value = "incorrect"
$ws.Cells.Item(22,3).Value = "correct"
$ws.Cells.Item(22,4).Value = 4.480341439799051
$ws.Cells.Item(22,5).Value = 2.776247428601131
$ws.Cells.Item(22,6).Value = 53.68394077873329
$ws.Cells.Item(22,7).Value = 8.266991960643239

$ws.Cells.Item(23,2).Value = "incorrect"
$ws.Cells.Item(23,3).Value = "correct"
$ws.Cells.Item(23,4).Value = 2.855144047978813
$ws.Cells.Item(23,5).Value = 4.55541127322908
$ws.Cells.Item(23,6).Value = 55.54942968569293
$ws.Cells.Item(23,7).Value = 74.77261233281929

$ws.Cells.Item(24,2).Value = "correct"
$ws.Cells.Item(24,3).Value = "incorrect"
$ws.Cells.Item(24,4).Value = 5.014938407765021
$ws.Cells.Item(24,5).Value = 1.765903471329896
$ws.Cells.Item(24,6).Value = 62.92587192860876
$ws.Cells.Item(24,7).Value = 28.64320358504051

$ws.Cells.Item(25,2).Value = "incorrect"
$ws.Cells.Item(25,3).Value = "correct"
$ws.Cells.Item(25,4).Value = 1.617373290308243
$ws.Cells.Item(25,5).Value = 3.236835757540429
$ws.Cells.Item(25,6).Value = 83.74947681738772
$ws.Cells.Item(25,7).Value = 20.37488766179312

$ws.Cells.Item(26,2).Value = "incorrect"
$ws.Cells.Item(26,3).Value = "correct"
$ws.Cells.Item(26,4).Value = 3.279412145024131
$ws.Cells.Item(26,5).Value = 2.748276710245083
$ws.Cells.Item(26,6).Value = 103.391944264572
$ws.Cells.Item(26,7).Value = 37.91628918599088

$ws.Cells.Item(27,2).Value = "incorrect"
$ws.Cells.Item(27,3).Value = "correct"
$ws.Cells.Item(27,4).Value = 1.268324829799188
$ws.Cells.Item(27,5).Value = 2.6881677829729
$ws.Cells.Item(27,6).Value = 96.42459495483435
$ws.Cells.Item(27,7).Value = 52.28600153926333

$ws.Cells.Item(28,2).Value = "incorrect"
$ws.Cells.Item(28,3).Value = "correct"
$ws.Cells.Item(28,4).Value = 3.305634670736914
$ws.Cells.Item(28,5).Value = 2.954636965382352
$ws.Cells.Item(28,6).Value = 103.8741821951283
$ws.Cells.Item(28,7).Value = 9.80815823909289

$ws.Cells.Item(29,2).Value = "correct"
$ws.Cells.Item(29,3).Value = "correct"
$ws.Cells.Item(29,4).Value = 4.11943106965656
$ws.Cells.Item(29,5).Value = 3.222156737950528
$ws.Cells.Item(29,6).Value = 44.98787656091629
$ws.Cells.Item(29,7).Value = 51.87109010646878

$ws.Cells.Item(30,2).Value = "incorrect"
$ws.Cells.Item(30,3).Value = "correct"
$ws.Cells.Item(30,4).Value = 1.46782996063464
$ws.Cells.Item(30,5).Value = 0.9354427891678414
$ws.Cells.Item(30,6).Value = 102.2973868545321
$ws.Cells.Item(30,7).Value = 16.99455144233779

$ws.Cells.Item(31,2).Value = "incorrect"
$ws.Cells.Item(31,3).Value = "correct"
$ws.Cells.Item(31,4).Value = 5.126430890697522
$ws.Cells.Item(31,5).Value = 2.103556019619417
$ws.Cells.Item(31,6).Value = 30.01080678860274
$ws.Cells.Item(31,7).Value = 8.278917356527998

$ws.Cells.Item(32,2).Value = "incorrect"
$ws.Cells.Item(32,3).Value = "correct"
$ws.Cells.Item(32,4).Value = 5.542401837991629
$ws.Cells.Item(32,5).Value = 1.770578489139599
$ws.Cells.Item(32,6).Value = 53.46239368871557
$ws.Cells.Item(32,7).Value = 25.64062600802293

$ws.Cells.Item(33,2).Value = "incorrect"
$ws.Cells.Item(33,3).Value = "incorrect"
$ws.Cells.Item(33,4).Value = 4.578020577209354
$ws.Cells.Item(33,5).Value = 4.439530085552135
$ws.Cells.Item(33,6).Value = 105.330048551436
$ws.Cells.Item(33,7).Value = 7.25977202303514

$ws.Cells.Item(34,2).Value = "incorrect"
$ws.Cells.Item(34,3).Value = "incorrect"
$ws.Cells.Item(34,4).Value = 5.536154577008677
$ws.Cells.Item(34,5).Value = 4.460778253903222
$ws.Cells.Item(34,6).Value = 46.94265933512254
$ws.Cells.Item(34,7).Value = 24.98682831039027

$ws.Cells.Item(35,2).Value = "incorrect"
$ws.Cells.Item(35,3).Value = "correct"
$ws.Cells.Item(35,4).Value = 4.070029369956201
$ws.Cells.Item(35,5).Value = 1.314548332149332
$ws.Cells.Item(35,6).Value = 99.98400643029271
$ws.Cells.Item(35,7).Value = 24.78546260633074

$ws.Cells.Item(36,2).Value = "incorrect"
$ws.Cells.Item(36,3).Value = "correct"
$ws.Cells.Item(36,4).Value = 4.13050602351068
$ws.Cells.Item(36,5).Value = 1.84923289102613
$ws.Cells.Item(36,6).Value = 38.18881980785324
$ws.Cells.Item(36,7).Value = 49.52627060765283

$ws.Cells.Item(37,2).Value = "correct"
$ws.Cells.Item(37,3).Value = "correct"
$ws.Cells.Item(37,4).Value = 3.740069442241901
$ws.Cells.Item(37,5).Value = 4.965834576475836
$ws.Cells.Item(37,6).Value = 69.75672752698343
$ws.Cells.Item(37,7).Value = 22.90176399753538

$ws.Cells.Item(38,2).Value = "correct"
$ws.Cells.Item(38,3).Value = "incorrect"
$ws.Cells.Item(38,4).Value = 2.024905419713451
$ws.Cells.Item(38,5).Value = 2.748428536197414
$ws.Cells.Item(38,6).Value = 107.8947483022908
$ws.Cells.Item(38,7).Value = 60.36424653746716

$ws.Cells.Item(39,2).Value = "incorrect"
$ws.Cells.Item(39,3).Value = "correct"
$ws.Cells.Item(39,4).Value = 3.364224436196156
$ws.Cells.Item(39,5).Value = 3.63289453373013
$ws.Cells.Item(39,6).Value = 109.0292371839198
$ws.Cells.Item(39,7).Value = 9.655572632712392

$ws.Cells.Item(40,2).Value = "incorrect"
$ws.Cells.Item(40,3).Value = "incorrect"
$ws.Cells.Item(40,4).Value = 1.00427808175715
$ws.Cells.Item(40,5).Value = 4.231066838955764
$ws.Cells.Item(40,6).Value = 27.10662553974002
$ws.Cells.Item(40,7).Value = 48.11671031583285

$ws.Cells.Item(41,2).Value = "correct"
$ws.Cells.Item(41,3).Value = "correct"
$ws.Cells.Item(41,4).Value = 2.327331747953587
$ws.Cells.Item(41,5).Value = 1.649648301856178
$ws.Cells.Item(41,6).Value = 83.02666032307613
$ws.Cells.Item(41,7).Value = 80.09307901952238

$ws.Cells.Item(42,2).Value = "incorrect"
$ws.Cells.Item(42,3).Value = "correct"
$ws.Cells.Item(42,4).Value = 3.586362490593167
$ws.Cells.Item(42,5).Value = 2.011741632832168
$ws.Cells.Item(42,6).Value = 40.23731518777672
$ws.Cells.Item(42,7).Value = 31.48549259107155

$ws.Cells.Item(43,2).Value = "incorrect"
$ws.Cells.Item(43,3).Value = "correct"
$ws.Cells.Item(43,4).Value = 5.324985203761658
$ws.Cells.Item(43,5).Value = 2.01208558039467
$ws.Cells.Item(43,6).Value = 91.02972873571527
$ws.Cells.Item(43,7).Value = 61.48018884849137

$ws.Cells.Item(44,2).Value = "incorrect"
$ws.Cells.Item(44,3).Value = "correct"
$ws.Cells.Item(44,4).Value = 5.580290780338933
$ws.Cells.Item(44,5).Value = 2.334300902026008
$ws.Cells.Item(44,6).Value = 53.20344050102381
$ws.Cells.Item(44,7).Value = 8.45553074898983

$ws.Cells.Item(45,2).Value = "correct"
$ws.Cells.Item(45,3).Value = "correct"
$ws.Cells.Item(45,4).Value = 2.639975240810683
$ws.Cells.Item(45,5).Value = 1.851722356976063
$ws.Cells.Item(45,6).Value = 10.19196290103863
$ws.Cells.Item(45,7).Value = 79.86685558489862

$ws.Cells.Item(46,2).Value = "correct"
$ws.Cells.Item(46,3).Value = "correct"
$ws.Cells.Item(46,4).Value = 5.083850765877295
$ws.Cells.Item(46,5).Value = 1.786542026895287
$ws.Cells.Item(46,6).Value = 70.31000835727554
$ws.Cells.Item(46,7).Value = 13.07178629525745

$ws.Cells.Item(47,2).Value = "incorrect"
$ws.Cells.Item(47,3).Value = "correct"
$ws.Cells.Item(47,4).Value = 2.475473898118651
$ws.Cells.Item(47,5).Value = 1.823435809622131
$ws.Cells.Item(47,6).Value = 31.37263717508207
$ws.Cells.Item(47,7).Value = 46.00664821327582

$ws.Cells.Item(48,2).Value = "incorrect"
$ws.Cells.Item(48,3).Value = "correct"
$ws.Cells.Item(48,4).Value = 5.112751282838721
$ws.Cells.Item(48,5).Value = 3.122664933071259
$ws.Cells.Item(48,6).Value = 77.55997068462763
$ws.Cells.Item(48,7).Value = 25.90278980941395

$ws.Cells.Item(49,2).Value = "incorrect"
$ws.Cells.Item(49,3).Value = "correct"
$ws.Cells.Item(49,4).Value = 5.043204928325432
$ws.Cells.Item(49,5).Value = 1.508193845843718
$ws.Cells.Item(49,6).Value = 98.79771630100844
$ws.Cells.Item(49,7).Value = 3.689468971030637

$ws.Cells.Item(50,2).Value = "incorrect"
$ws.Cells.Item(50,3).Value = "incorrect"
$ws.Cells.Item(50,4).Value = 3.242160495395595
$ws.Cells.Item(50,5).Value = 3.425095578988311
$ws.Cells.Item(50,6).Value = 52.27609269247633
$ws.Cells.Item(50,7).Value = 31.07932218109558

$ws.Cells.Item(51,2).Value = "incorrect"
$ws.Cells.Item(51,3).Value = "correct"
$ws.Cells.Item(51,4).Value = 3.582010385333696
$ws.Cells.Item(51,5).Value = 2.633178368107929
$ws.Cells.Item(51,6).Value = 88.51801641846555
$ws.Cells.Item(51,7).Value = 37.75881390453094
